$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the duty time shared strings: remove spaces around the dash,
# and fix the second shift's end time (10am -> 10pm).
$ws.Range("A1").Value = "6am-2pm"
$ws.Range("A2").Value = "2pm-10pm"

# Update the active selection on the sheet view.
$ws.Range("B4").Select()
